# major accuracy check update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the polyAIsolationProtocol value (column G) for every data row:
# "NEBNextPoly(A)E7490" -> "NEBNextPoly(A)E7490L"
$ws.Range("G2:G41").Value = "NEBNextPoly(A)E7490L"

# Widen column G to fit the longer label (~35.98 characters).
$ws.Columns.Item(7).ColumnWidth = 35.17

# Restore the view: scroll back to the top-left and select the column
# that was just reviewed/updated.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G2:G41").Select()
